# 1. Update Metadata sheet: Total Queries count 3 -> 4
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = 4

# 2. Add new worksheet "distance from Dma50" as the last sheet (after "1 Month Performance")
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "distance from Dma50"

# 3. Header row (bold, centered, thin border - matches style used on other sheets)
$ws.Cells.Item(1,1).Value = "Icon"
$ws.Cells.Item(1,2).Value = "Stock"
$ws.Cells.Item(1,3).Value = "Distance From Sma50"
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# 4. Data rows
$ws.Cells.Item(2,1).Value = "📈"
$ws.Cells.Item(2,2).Value = "NIFTYPSUBANK"
$ws.Cells.Item(2,3).Value = 10.196
$ws.Cells.Item(3,1).Value = "📈"
$ws.Cells.Item(3,2).Value = "NIFTYMETAL"
$ws.Cells.Item(3,3).Value = 8.6656
$ws.Cells.Item(4,1).Value = "📈"
$ws.Cells.Item(4,2).Value = "NIFTYOILANDGAS"
$ws.Cells.Item(4,3).Value = 6.5062
$ws.Cells.Item(5,1).Value = "📈"
$ws.Cells.Item(5,2).Value = "CNXREALTY"
$ws.Cells.Item(5,3).Value = 5.7914
$ws.Cells.Item(6,1).Value = "📈"
$ws.Cells.Item(6,2).Value = "NIFTYCOMMODITIES"
$ws.Cells.Item(6,3).Value = 5.7904
$ws.Cells.Item(7,1).Value = "📈"
$ws.Cells.Item(7,2).Value = "CNXINFRA"
$ws.Cells.Item(7,3).Value = 5.6788
$ws.Cells.Item(8,1).Value = "📈"
$ws.Cells.Item(8,2).Value = "NIFTYPVTBANK"
$ws.Cells.Item(8,3).Value = 5.0334
$ws.Cells.Item(9,1).Value = "📈"
$ws.Cells.Item(9,2).Value = "BANKNIFTY"
$ws.Cells.Item(9,3).Value = 4.9435
$ws.Cells.Item(10,1).Value = "📈"
$ws.Cells.Item(10,2).Value = "NIFTYFINSERVICE"
$ws.Cells.Item(10,3).Value = 4.0456
$ws.Cells.Item(11,1).Value = "📈"
$ws.Cells.Item(11,2).Value = "NIFTYMIDCAP50"
$ws.Cells.Item(11,3).Value = 3.993
$ws.Cells.Item(12,1).Value = "📈"
$ws.Cells.Item(12,2).Value = "CNXENERGY"
$ws.Cells.Item(12,3).Value = 3.888
$ws.Cells.Item(13,1).Value = "📈"
$ws.Cells.Item(13,2).Value = "NIFTY"
$ws.Cells.Item(13,3).Value = 3.7752
$ws.Cells.Item(14,1).Value = "📈"
$ws.Cells.Item(14,2).Value = "CNXMIDCAP"
$ws.Cells.Item(14,3).Value = 3.6938
$ws.Cells.Item(15,1).Value = "📈"
$ws.Cells.Item(15,2).Value = "NIFTY200"
$ws.Cells.Item(15,3).Value = 3.6551
$ws.Cells.Item(16,1).Value = "📈"
$ws.Cells.Item(16,2).Value = "NIFTY100"
$ws.Cells.Item(16,3).Value = 3.6395
$ws.Cells.Item(17,1).Value = "📈"
$ws.Cells.Item(17,2).Value = "NIFTY500"
$ws.Cells.Item(17,3).Value = 3.3949
$ws.Cells.Item(18,1).Value = "📈"
$ws.Cells.Item(18,2).Value = "NIFTYCPSE"
$ws.Cells.Item(18,3).Value = 3.0412
$ws.Cells.Item(19,1).Value = "📈"
$ws.Cells.Item(19,2).Value = "NIFTY50VALUE20"
$ws.Cells.Item(19,3).Value = 2.9668
$ws.Cells.Item(20,1).Value = "📈"
$ws.Cells.Item(20,2).Value = "CNXSMALLCAP"
$ws.Cells.Item(20,3).Value = 2.965
$ws.Cells.Item(21,1).Value = "📈"
$ws.Cells.Item(21,2).Value = "CNXNIFTYJUNIOR"
$ws.Cells.Item(21,3).Value = 2.9217
$ws.Cells.Item(22,1).Value = "📈"
$ws.Cells.Item(22,2).Value = "CNXIT"
$ws.Cells.Item(22,3).Value = 2.1425
$ws.Cells.Item(23,1).Value = "📈"
$ws.Cells.Item(23,2).Value = "NIFTYHEALTHCARE"
$ws.Cells.Item(23,3).Value = 2.1379
$ws.Cells.Item(24,1).Value = "📈"
$ws.Cells.Item(24,2).Value = "NIFTYCONSUMPTION"
$ws.Cells.Item(24,3).Value = 2.092
$ws.Cells.Item(25,1).Value = "📈"
$ws.Cells.Item(25,2).Value = "NIFTYAUTO"
$ws.Cells.Item(25,3).Value = 1.6784
$ws.Cells.Item(26,1).Value = "📈"
$ws.Cells.Item(26,2).Value = "NIFTYGROWSECT15"
$ws.Cells.Item(26,3).Value = 1.5497
$ws.Cells.Item(27,1).Value = "📈"
$ws.Cells.Item(27,2).Value = "CNXPHARMA"
$ws.Cells.Item(27,3).Value = 1.546
$ws.Cells.Item(28,1).Value = "📈"
$ws.Cells.Item(28,2).Value = "NIFTYFMCG"
$ws.Cells.Item(28,3).Value = 1.3946
$ws.Cells.Item(29,1).Value = "📈"
$ws.Cells.Item(29,2).Value = "NIFTYCONSURDURBL"
$ws.Cells.Item(29,3).Value = 0.3549
$ws.Cells.Item(30,1).Value = "📈"
$ws.Cells.Item(30,2).Value = "NIFTYMEDIA"
$ws.Cells.Item(30,3).Value = -2.0693
$ws.Cells.Item(31,1).Value = "📈"
$ws.Cells.Item(31,2).Value = "NIFTYFINSEREXBNK"
$ws.Cells.Item(31,3).Value = "N/A"
$ws.Cells.Item(32,1).Value = "📈"
$ws.Cells.Item(32,2).Value = "NIFTYMSITTELCM"
$ws.Cells.Item(32,3).Value = "N/A"
$ws.Cells.Item(33,1).Value = "📈"
$ws.Cells.Item(33,2).Value = "NIFTYMSFINSERV"
$ws.Cells.Item(33,3).Value = "N/A"
